$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (includes two rank swaps: Aave/MXToken rows 40-41, Quant/PaxDollar rows 49-50)
# Leading apostrophe forces text interpretation so numeric-looking strings (e.g. "1.00", "34.447.50") are
# preserved exactly as text instead of being auto-converted to numbers; Style reset avoids introducing a new
# "Text" number-format style so the cell keeps its original (default) style.
$ws.Range("D2").Value = "'34.447.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.40%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.803.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.97%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.30%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'227.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.15%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.575"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.50%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value = "'36.27"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +4.21%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.58%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D11").Value = "'0.0964"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.38%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.063.62"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.03%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'11.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.21%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'1.806.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.82%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -0.33%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +3.18%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'34.427.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.45%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'70.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.10%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'244.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.22%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0789"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'11.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.21%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.30%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.37%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'173.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.95%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +3.65%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +7.77%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.07%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +1.00%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.23%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -0.76%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.08%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.60%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.56%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -2.70%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.392.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.79%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.672"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.19%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -6.80%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.89%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -1.31%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Aave"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'82.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -4.88%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'MXToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'2.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.66%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +0.01%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.47%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +8.21%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'13.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'6.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.96%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0502"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -4.21%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.964.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.10%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Quant"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'104.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.71%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'PaxDollar"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.32%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.62%  "
$ws.Range("E51").Style = "Normal"

Write-Host "Applied 81 cell updates"
